$wb = $excel.ActiveWorkbook

# --- Update Layer0 sheet (sheet1) ---
$ws1 = $wb.Worksheets.Item("Layer0")

$ws1.Range("B2").Value = -96.46666825684912
$ws1.Range("C2").Value = -0.09480787604336607
$ws1.Range("D2").Value = -100.4919959065659
$ws1.Range("E2").Value = 60.37763814787296
$ws1.Range("F2").Value = -156.9944966294624
$ws1.Range("G2").Value = 59.43286883802325
$ws1.Range("H2").Value = -165.0627348017731
$ws1.Range("I2").Value = 46.00518855658557
$ws1.Range("J2").Value = -160.5746313205443
$ws1.Range("K2").Value = -150.8736514314649
$ws1.Range("L2").Value = -134.8530198390092
$ws1.Range("M2").Value = -22.15070272920044
$ws1.Range("N2").Value = 98.03905139260361
$ws1.Range("O2").Value = 60.17265246174674
$ws1.Range("P2").Value = -106.8984912012236
$ws1.Range("Q2").Value = -28.25289488435505

$ws1.Range("B3").Value = 87.85082751276639
$ws1.Range("C3").Value = 116.6434496303378
$ws1.Range("D3").Value = 86.001863293495
$ws1.Range("E3").Value = -65.74184547215084
$ws1.Range("F3").Value = 40.93851092644854
$ws1.Range("G3").Value = 90.97078417504368
$ws1.Range("H3").Value = 60.1083072012806
$ws1.Range("I3").Value = 41.69949760054234
$ws1.Range("J3").Value = 53.56230620112576
$ws1.Range("K3").Value = 57.96677939056418
$ws1.Range("L3").Value = 72.51437418969205
$ws1.Range("M3").Value = 85.44391387466088
$ws1.Range("N3").Value = 79.69638968636215
$ws1.Range("O3").Value = 96.55222620087987
$ws1.Range("P3").Value = 86.76760605967978
$ws1.Range("Q3").Value = 92.24386123399943

$ws1.Range("B4").Value = 91.45474791632861
$ws1.Range("C4").Value = 48.51969670222189
$ws1.Range("D4").Value = 33.49191840994484
$ws1.Range("E4").Value = 53.34027552622367
$ws1.Range("F4").Value = 42.24285713075309
$ws1.Range("G4").Value = -26.95708167229242
$ws1.Range("H4").Value = 47.65833943566464
$ws1.Range("I4").Value = 87.81819203555764
$ws1.Range("J4").Value = 40.83401713897556
$ws1.Range("K4").Value = 73.40654448479631
$ws1.Range("L4").Value = 49.03276302486768
$ws1.Range("M4").Value = 95.34393966069389
$ws1.Range("N4").Value = -5.85298353227508
$ws1.Range("O4").Value = 23.2018266029816
$ws1.Range("P4").Value = 54.92520732930956
$ws1.Range("Q4").Value = 19.84012212555737

$ws1.Range("B5").Value = 10.41030422579971
$ws1.Range("C5").Value = 11.53046999478885
$ws1.Range("D5").Value = 25.36204848728877
$ws1.Range("E5").Value = 8.867389129383875
$ws1.Range("F5").Value = 16.61781095801848
$ws1.Range("G5").Value = 7.8759674249787
$ws1.Range("H5").Value = 30.69137647333278
$ws1.Range("I5").Value = 0.7121046719257649
$ws1.Range("J5").Value = 35.98398155076277
$ws1.Range("K5").Value = 15.66387478528744
$ws1.Range("L5").Value = 34.34778473193298
$ws1.Range("M5").Value = 10.87396035519473
$ws1.Range("N5").Value = 8.861756689302434
$ws1.Range("O5").Value = 22.8127018823403
$ws1.Range("P5").Value = 46.13364681572625
$ws1.Range("Q5").Value = 70.86030672892906

$ws1.Range("B6").Value = 15.53388870076511
$ws1.Range("C6").Value = 6.700080489138944
$ws1.Range("D6").Value = 42.74571363457571
$ws1.Range("E6").Value = 14.71931060869734
$ws1.Range("F6").Value = 17.04181719727595
$ws1.Range("G6").Value = 24.49043943533821
$ws1.Range("H6").Value = 24.139550078603
$ws1.Range("I6").Value = 7.013120278892869
$ws1.Range("J6").Value = 18.09638600167139
$ws1.Range("K6").Value = 1.857394809978736
$ws1.Range("L6").Value = 15.68051131621504
$ws1.Range("M6").Value = 13.06876856563135
$ws1.Range("N6").Value = 2.962069356858701
$ws1.Range("O6").Value = 0.9340544062732621
$ws1.Range("P6").Value = 9.983851734927446
$ws1.Range("Q6").Value = 13.38803031540529

$ws1.Range("B7").Value = 15.21032549298393
$ws1.Range("C7").Value = 7.728995083865871
$ws1.Range("D7").Value = -1.552963196010012
$ws1.Range("E7").Value = 6.42490974394512
$ws1.Range("F7").Value = -0.1843261278838499
$ws1.Range("G7").Value = -3.641569767994171
$ws1.Range("H7").Value = 2.576808402249654
$ws1.Range("I7").Value = 5.763245046133774
$ws1.Range("J7").Value = -2.432806300554565
$ws1.Range("K7").Value = 12.31364635709139
$ws1.Range("L7").Value = 0.6190783370522572
$ws1.Range("M7").Value = 4.411664937249305
$ws1.Range("N7").Value = -0.7064248728435306
$ws1.Range("O7").Value = 5.8737102799694
$ws1.Range("P7").Value = 1.276548924858215
$ws1.Range("Q7").Value = 0.1040451786725431

$ws1.Range("B8").Value = -3.006573661561808
$ws1.Range("C8").Value = 1.897718660232665
$ws1.Range("D8").Value = -2.065122003052232
$ws1.Range("E8").Value = 5.44227760587315
$ws1.Range("F8").Value = 2.853661471441971
$ws1.Range("G8").Value = -4.737158062427254
$ws1.Range("H8").Value = -4.648281149196145
$ws1.Range("I8").Value = 5.181526222559645
$ws1.Range("J8").Value = -0.7160103480767851
$ws1.Range("K8").Value = -2.45921628038894
$ws1.Range("L8").Value = 3.134907327161828
$ws1.Range("M8").Value = -0.1143434045131765
$ws1.Range("N8").Value = -0.6716865449645854
$ws1.Range("O8").Value = 3.481826853836204
$ws1.Range("P8").Value = 8.341370802047996
$ws1.Range("Q8").Value = 8.529465923360503

$ws1.Range("B9").Value = -1.963698006972093
$ws1.Range("C9").Value = 4.658382246481182
$ws1.Range("D9").Value = 16.98674377688561
$ws1.Range("E9").Value = -0.1303112556414516
$ws1.Range("F9").Value = 3.633433547363901
$ws1.Range("G9").Value = 0.8676011996169525
$ws1.Range("H9").Value = 1.806921842145263
$ws1.Range("I9").Value = -0.1289539537934635
$ws1.Range("J9").Value = 5.807773347432053
$ws1.Range("K9").Value = 3.431812179367124
$ws1.Range("L9").Value = 3.070574566612545
$ws1.Range("M9").Value = -5.541300647940449
$ws1.Range("N9").Value = 1.514618187954278
$ws1.Range("O9").Value = 4.470482461895665
$ws1.Range("P9").Value = -0.6728589738703831
$ws1.Range("Q9").Value = 3.676153684553503

# --- Update Layer1 sheet (sheet2) ---
$ws2 = $wb.Worksheets.Item("Layer1")

$ws2.Range("B2").Value = -120.96144736032
$ws2.Range("B3").Value = 203.5245476012108
$ws2.Range("B4").Value = 175.8353032697162
$ws2.Range("B5").Value = 189.3482587257351
$ws2.Range("B6").Value = 204.2174265515936
$ws2.Range("B7").Value = 252.3087708012458
$ws2.Range("B8").Value = 171.6691488211002
$ws2.Range("B9").Value = 280.5750348098936
$ws2.Range("B10").Value = 191.2706505513266
$ws2.Range("B11").Value = 270.5675106169301
$ws2.Range("B12").Value = 235.0905588785757
$ws2.Range("B13").Value = 244.6378895374059
$ws2.Range("B14").Value = 205.1207063958388
$ws2.Range("B15").Value = 198.7650628732942
$ws2.Range("B16").Value = 183.4333072697281
$ws2.Range("B17").Value = 221.5530269450389
$ws2.Range("B18").Value = 173.5123285457267
